$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contract Parts")

# AG8:AG12 used to hold a static serial-date value (45222); the author
# resolved a conflict by swapping these in for a volatile =TODAY() formula.
# Write each cell individually so the formula isn't exported as a single
# shared-formula group.
$ws.Range("AG8").Formula  = "=TODAY()"
$ws.Range("AG9").Formula  = "=TODAY()"
$ws.Range("AG10").Formula = "=TODAY()"
$ws.Range("AG11").Formula = "=TODAY()"
$ws.Range("AG12").Formula = "=TODAY()"

# Move the saved selection/active cell from AG15 to AG9.
$null = $ws.Range("AG9").Select()
